$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing ParentID ("2") for the section-header rows that were
# previously left blank in column B.
$ws.Range("B12").Value = "2"
$ws.Range("B15").Value = "2"
$ws.Range("B18").Value = "2"
$ws.Range("B20").Value = "2"

# Restore the saved cursor/selection position to B27.
$ws.Range("B27").Select() | Out-Null
